$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old weekly date columns (B:G). Excel shifts the remaining
# "11_02_2024" / "18_02_2024" columns (old H:I) left into B:C.
$ws.Range("B:G").Delete()

# Re-apply the underline style that used to live on the old D6/F6 cells
# (style index 1) so those positions keep their formatting even though
# their values are gone, matching the reorganized layout.
$ws.Range("D6").Font.Underline = $true
$ws.Range("F6").Font.Underline = $true

$ws.Range("I11").Select() | Out-Null
